# Deploy the implementation guide:
#  - bump the "Date" metadata value to the new publish timestamp
#  - add a new "EXOR / Exomiser Report" concept row to the Concepts sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value (B8) ---------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-02-20T18:51:53+00:00"

# --- Concepts sheet: append a new concept row (row 8) ----------------------
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Copy the last data row (row 7) down to the new row so the new row inherits
# the same cell formatting/style (borders, alignment, etc.) as the rest of
# the table, including keeping column A as a shared-string "1" rather than
# a numeric literal.
$wsConcepts.Range("A7:D7").Copy()
$wsConcepts.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats

# Re-copy just A7 (Level = "1") over A8 so the cell keeps its original
# text type/value instead of being left blank by the formats-only paste.
$wsConcepts.Range("A7").Copy()
$wsConcepts.Range("A8").PasteSpecial(-4163)  # xlPasteValues

# Fill in the new concept's Code and Display.
$wsConcepts.Range("B8").Value = "EXOR"
$wsConcepts.Range("C8").Value = "Exomiser Report"
